$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.547.74"
$ws.Range("E2").Value = "  +0.83%  "

$ws.Range("D3").Value = "3.099.19"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.090.51"
$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.05%  "

$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").Value = "3.611.30"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").Value = "63.358.37"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.095.07"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("E22").Value = "  -0.63%  "

$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.88%  "

$ws.Range("E29").Value = "  +0.06%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "0.0₃0849"
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.31%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "434.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").Value = "2.880.59"
$ws.Range("E44").Value = "  -1.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "

$ws.Range("E46").Value = "  -2.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.75%  "

